$d = $word.ActiveDocument
$rng1 = $d.Content
$rng1.Find.Execute("Może dodawać nowe utwory,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p1 = $rng1.Paragraphs(1)
$p1End = $p1.Range
$ins1 = $d.Range($p1End.End - 1, $p1End.End - 1)
$ins1.InsertAfter(" //nie udało się wprowadzić")

$rng2 = $d.Content
$rng2.Find.Execute("Może nadawać uprawnienia administratora;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2 = $rng2.Paragraphs(1)
$p2End = $p2.Range
$ins2 = $d.Range($p2End.End - 1, $p2End.End - 1)
$ins2.InsertAfter(" //nie udało się wprowadzić")

$paras = $d.Paragraphs
$last = $paras.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$paras2 = $d.Paragraphs
$newLast = $paras2.Last
$newLast.Range.Text = "JS nie działa."
